$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D and E columns are treated as text so numeric-looking values are not
# auto-converted to numbers by Excel (matches original inlineStr text cells).
$rngDE = $ws.Range("D2:E51")
$rngDE.NumberFormat = "@"

$ws.Range("D2").Value = "29.404.56"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "1.877.48"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "0.7177"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("D6").Value = "243.74"
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "0.07968"
$ws.Range("E8").Value = "  +1.19%  "
$ws.Range("D9").Value = "0.3150"
$ws.Range("D10").Value = "24.98"
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("D11").Value = "0.08128"
$ws.Range("E11").Value = "  -3.10%  "
$ws.Range("D12").Value = "1.881.48"
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("D13").Value = "95.20"
$ws.Range("E13").Value = "  +4.13%  "
$ws.Range("E14").Value = "  -0.18%  "
$ws.Range("D15").Value = "0.7076"
$ws.Range("E15").Value = "  -1.55%  "
$ws.Range("D16").Value = "6.411"
$ws.Range("E16").Value = "  +4.22%  "
$ws.Range("D17").Value = "0.000008447"
$ws.Range("E17").Value = "  +0.99%  "
$ws.Range("D18").Value = "29.405.96"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "253.56"
$ws.Range("E19").Value = "  +5.20%  "
$ws.Range("D20").Value = "13.39"
$ws.Range("E20").Value = "  +1.12%  "
$ws.Range("D21").Value = "2.136.45"
$ws.Range("E21").Value = "  +0.11%  "
$ws.Range("E22").Value = "  -0.18%  "
$ws.Range("D23").Value = "7.671"
$ws.Range("E23").Value = "  -1.57%  "
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").Value = "0.1590"
$ws.Range("E25").Value = "  -0.70%  "
$ws.Range("D26").Value = "9.066"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("D27").Value = "161.98"
$ws.Range("D28").Value = "18.92"
$ws.Range("E28").Value = "  +1.95%  "
$ws.Range("D29").Value = "1.508"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").Value = "4.305"
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("D32").Value = "1.218"
$ws.Range("E32").Value = "  -1.01%  "
$ws.Range("D33").Value = "0.05334"
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("D34").Value = "1.950"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").Value = "0.7566"
$ws.Range("E35").Value = "  +1.16%  "
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("D37").Value = "2.702"
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("D38").Value = "0.01896"
$ws.Range("E38").Value = "  +0.64%  "
$ws.Range("D39").Value = "1.272.81"
$ws.Range("E39").Value = "  -2.21%  "
$ws.Range("E40").Value = "  +0.98%  "
$ws.Range("D41").Value = "6.391"
$ws.Range("E41").Value = "  -2.00%  "
$ws.Range("D44").Value = "74.36"
$ws.Range("E44").Value = "  +2.28%  "
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("E46").Value = "  -2.86%  "
$ws.Range("D47").Value = "2.032.26"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("D49").Value = "0.5198"
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("D50").Value = "9.528"
$ws.Range("E50").Value = "  +0.69%  "
$ws.Range("E51").Value = "  -0.32%  "

# Row 42/43: TrustWalletToken and Quant swap positions
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "112.12"
$ws.Range("E42").Value = "  +1.29%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "0.9061"
$ws.Range("E43").Value = "  +1.56%  "

# Restore default (unstyled) style for the D:E range so only the number format
# change used to force text entry is undone, matching original formatting.
$rngDE.Style = "Normal"
